# Weekly update: insert a new week's pricing data (rows 182-183) for
# "Terminal La Palmera de La Serena - Betarraga" and push the existing
# rows down, growing the used range from A1:R185 to A1:R187.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 182:183 - this shifts the former rows
# 182-185 down to 184-187 (Excel copies formatting, incl. the date
# style on column D, from the row above when inserting).
$ws.Rows("182:183").Insert()

# Fill in the new week's data (D = 44509, the newest date) in the two
# freshly inserted rows.
$ws.Range("A182").Value = 8
$ws.Range("B182").Value = "Terminal La Palmera de La Serena"
$ws.Range("C182").Value = "Coquimbo"
$ws.Range("D182").Value = 44509
$ws.Range("E182").Value = 4
$ws.Range("F182").Value = 100114014
$ws.Range("G182").Value = "Betarraga"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 2940
$ws.Range("K182").Value = 450
$ws.Range("L182").Value = 500
$ws.Range("M182").Value = 475
$ws.Range("N182").Value = "`$/paquete 3 unidades"
$ws.Range("O182").Value = "Provincia del Elquí"
$ws.Range("P182").Value = 158
$ws.Range("Q182").Value = 3
$ws.Range("R182").Value = "Hortaliza"

$ws.Range("A183").Value = 8
$ws.Range("B183").Value = "Terminal La Palmera de La Serena"
$ws.Range("C183").Value = "Coquimbo"
$ws.Range("D183").Value = 44509
$ws.Range("E183").Value = 4
$ws.Range("F183").Value = 100114014
$ws.Range("G183").Value = "Betarraga"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Segunda"
$ws.Range("J183").Value = 1500
$ws.Range("K183").Value = 350
$ws.Range("L183").Value = 400
$ws.Range("M183").Value = 375
$ws.Range("N183").Value = "`$/paquete 3 unidades"
$ws.Range("O183").Value = "Provincia del Elquí"
$ws.Range("P183").Value = 125
$ws.Range("Q183").Value = 3
$ws.Range("R183").Value = "Hortaliza"

Write-Host "Done inserting and filling rows 182-183"
